$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 17 de Agosto de 2020 a las 12:41"

# Full refreshed country table (rank-sorted by "Casos totales" desc), rows 4-219
$data = @(
    @("Estados Unidos",5567765,1133,2922929,2471697,0,11,173139),
    @("Brasil",3340197,0,2432456,799862,0,0,107879),
    @("India",2651290,3974,1920265,679946,0,34,51079),
    @("Rusia",927745,4892,736101,175904,0,55,15740),
    @("Sudafrica",587345,0,472377,103129,0,0,11839),
    @("Peru",535946,0,365367,144298,0,0,26281),
    @("Mexico",522162,4448,355101,110304,0,214,56757),
    @("Colombia",468332,0,287436,165799,0,0,15097),
    @("Chile",385946,0,358828,16666,0,0,10452),
    @("España",358843,0,0,0,0,0,28617),
    @("Iran",345450,2247,299157,26489,0,165,19804),
    @("Reino Unido",318484,0,0,0,0,0,41366),
    @("Arabia Saudita",298542,0,266953,28181,0,0,3408),
    @("Argentina",294569,0,211702,77164,0,0,5703),
    @("Pakistan",289215,498,269087,13953,0,7,6175),
    @("Banglades",279144,2595,160591,114859,0,37,3694),
    @("Italia",253915,0,203786,14733,0,0,35396),
    @("Turquia",249309,0,230969,12366,0,0,5974),
    @("Alemania",224997,0,202900,12807,0,0,9290),
    @("Francia",218536,0,83848,104278,0,0,30410),
    @("Irak",176931,0,125374,45697,0,0,5860),
    @("Filipinas",164474,3314,112759,49034,0,18,2681),
    @("Indonesia",141370,1821,94458,40705,0,57,6207),
    @("Canada",122087,0,108484,4577,0,0,9026),
    @("Catar",115080,0,111794,3093,0,0,193),
    @("Kazajistan",103033,337,82777,18987,0,0,1269),
    @("Ecuador",101542,0,87022,8450,0,0,6070),
    @("Bolivia",100344,1198,36491,59795,0,55,4058),
    @("Egipto",96475,0,59743,31572,0,0,5160),
    @("Israel",93691,1011,69062,23939,0,5,690),
    @("Ucrania",92820,1464,48164,42567,0,19,2089),
    @("Republica Dominicana",86309,0,52905,31951,0,0,1453),
    @("China",84849,22,79603,612,0,0,4634),
    @("Suecia",84294,0,0,0,0,0,5783),
    @("Oman",83226,140,77812,4826,0,16,588),
    @("Panama",81940,0,55001,25172,0,0,1767),
    @("Belgica",78323,454,17994,50390,0,4,9939),
    @("Kuwait",76205,0,68135,7569,0,0,501),
    @("Rumania",71194,733,32759,35406,0,38,3029),
    @("Bielorrusia",69516,0,67072,1834,0,0,610),
    @("Emiratos Arabes Unidos",64312,0,57694,6254,0,0,364),
    @("Paises Bajos",63002,0,0,0,0,0,6172),
    @("Guatemala",62562,0,50692,9491,0,0,2379),
    @("Polonia",57279,595,39359,16035,0,8,1885),
    @("Singapur",55838,91,51953,3858,0,0,27),
    @("Japon",54714,0,40080,13546,0,0,1088),
    @("Portugal",54102,0,39697,12627,0,0,1778),
    @("Honduras",50502,523,7339,41588,0,8,1575),
    @("Nigeria",49068,0,36497,11596,0,0,975),
    @("Barein",46835,0,43128,3535,0,2,172),
    @("Ghana",42653,121,40567,1847,0,8,239),
    @("Marruecos",42489,0,29344,12487,0,0,658),
    @("Kirguistan",41991,135,34537,5958,0,1,1496),
    @("Armenia",41701,38,34655,6222,0,6,824),
    @("Argelia",38583,0,27017,10196,0,0,1370),
    @("Suiza",38124,0,33300,2833,0,0,1991),
    @("Afganistan",37596,0,27166,9055,0,0,1375),
    @("Uzbekistan",35513,184,30973,4306,0,2,234),
    @("Azerbaiyan",34219,0,31875,1838,0,0,506),
    @("Venezuela",33755,0,22700,10774,0,0,281),
    @("Moldavia",30183,0,21220,8067,0,0,896),
    @("Kenia",30120,0,16656,12990,0,0,474),
    @("Etiopia",29876,0,12359,16989,0,0,528),
    @("Serbia",29682,0,27061,1947,0,0,674),
    @("Costa Rica",28465,0,9062,19109,0,0,294),
    @("Irlanda",27257,0,23364,2119,0,0,1774),
    @("Nepal",26660,0,17335,9221,0,0,104),
    @("Australia",23558,270,14080,9057,0,25,421),
    @("Austria",23534,164,20765,2040,0,1,729),
    @("El Salvador",22912,0,10814,11480,0,6,618),
    @("Chequia",20012,0,13799,5816,0,0,397),
    @("Camerun",18469,0,16540,1528,0,0,401),
    @("Costa de Marfil",17026,0,13947,2969,0,0,110),
    @("Estado de Palestina",16534,0,9838,6586,0,0,110),
    @("Bosnia y Herzegovina",15801,0,9619,5711,0,0,471),
    @("Dinamarca",15617,0,13340,1656,0,0,621),
    @("Corea del Sur",15515,197,13917,1293,0,0,305),
    @("Bulgaria",14365,0,9186,4681,0,0,498),
    @("Madagascar",13886,59,12603,1112,0,1,171),
    @("Republica de Macedonia",12739,0,9174,3021,0,0,544),
    @("Sudan",12314,0,6350,5166,0,0,798),
    @("Senegal",12162,0,7677,4232,0,0,253),
    @("Noruega",10005,0,8857,887,0,0,261),
    @("Paraguay",9791,0,6034,3619,0,0,138),
    @("Consejo Danes para los Refugiados",9676,0,8705,731,0,0,240),
    @("Zambia",9343,0,8412,671,0,0,260),
    @("Malasia",9212,12,8876,211,0,0,125),
    @("Libano",8881,0,2724,6054,0,0,103),
    @("Guayana Francesa",8588,0,7893,642,0,0,53),
    @("Guinea",8482,0,7364,1067,0,0,51),
    @("Gabon",8225,0,6277,1897,0,0,51),
    @("Libia",8172,0,933,7086,0,0,153),
    @("Tayikistan",8065,0,6855,1146,0,0,64),
    @("Haiti",7879,0,5235,2448,0,0,196),
    @("Finlandia",7752,21,7050,369,0,0,333),
    @("Luxemburgo",7458,0,6500,835,0,0,123),
    @("Albania",7380,0,3794,3358,0,0,228),
    @("Grecia",7075,0,3804,3043,0,0,228),
    @("Mauritania",6701,0,5985,559,0,0,157),
    @("Croacia",6571,0,5220,1185,0,0,166),
    @("Maldivas",5785,0,3349,2414,0,0,22),
    @("Republica de Yibuti",5369,0,5202,108,0,0,59),
    @("Zimbabue",5261,0,2092,3037,0,0,132),
    @("Malaui",5072,0,2626,2285,0,0,161),
    @("Hungria",4946,30,3630,708,0,0,608),
    @("Guinea Ecuatorial",4821,0,2182,2556,0,0,83),
    @("Republica de Africa Central",4652,0,1728,2863,0,0,61),
    @("Hong Kong",4525,44,3599,857,0,0,69),
    @("Namibia",4154,0,2370,1749,0,0,35),
    @("Nicaragua",4115,0,2913,1074,0,0,128),
    @("Montenegro",4035,0,2910,1048,0,0,77),
    @("Suazilandia",3839,0,2268,1501,0,0,70),
    @("Congo",3831,0,1625,2130,0,0,76),
    @("Tailandia",3378,1,3194,126,0,0,58),
    @("Cuba",3316,0,2620,608,0,0,88),
    @("Somalia",3256,0,2374,789,0,0,93),
    @("Cabo Verde",3179,0,2317,827,0,0,35),
    @("Mayotte",3119,0,2940,140,0,0,39),
    @("Surinam",3016,0,2036,933,0,0,47),
    @("Eslovaquia",2907,5,1969,907,0,0,31),
    @("Sri Lanka",2893,0,2676,206,0,0,11),
    @("Mozambique",2855,0,1163,1673,0,0,19),
    @("Mali",2640,0,1987,528,0,0,125),
    @("Sudan del Sur",2489,0,1175,1267,0,0,47),
    @("Ruanda",2453,0,1648,797,0,0,8),
    @("Lituania",2436,20,1705,650,0,0,81),
    @("Eslovenia",2429,13,2052,248,0,0,129),
    @("Estonia",2192,2,1976,153,0,0,63),
    @("Guinea-Bisau",2117,0,1015,1069,0,0,33),
    @("Tunez",2107,0,1358,695,0,0,54),
    @("Benin",2063,0,1690,334,0,0,39),
    @("Islandia",2011,0,1880,121,0,0,10),
    @("Sierra Leona",1956,0,1506,381,0,0,69),
    @("Angola",1906,0,628,1190,0,0,88),
    @("Gambia",1872,0,401,1408,0,0,63),
    @("Yemen",1869,0,1013,326,0,0,530),
    @("Siria",1677,0,417,1196,0,0,64),
    @("Nueva Zelanda",1631,9,1531,78,0,0,22),
    @("Uganda",1500,0,1142,345,0,0,13),
    @("Uruguay",1440,0,1200,202,0,0,38),
    @("Jordania",1378,0,1236,131,0,0,11),
    @("Malta",1375,69,759,607,0,0,9),
    @("Georgia",1341,5,1092,232,0,0,17),
    @("Republica de Chipre",1339,0,870,449,0,0,20),
    @("Letonia",1323,1,1078,213,0,0,32),
    @("Bahamas",1315,0,189,1108,0,0,18),
    @("Burkina Faso",1267,0,1013,199,0,0,55),
    @("Liberia",1257,0,788,387,0,0,82),
    @("Botsuana",1214,0,120,1091,0,0,3),
    @("Niger",1167,0,1078,20,0,0,69),
    @("Togo",1147,0,843,277,0,0,27),
    @("Jamaica",1113,7,764,335,0,0,14),
    @("Aruba",1102,0,200,898,0,0,4),
    @("Principado de Andorra",989,0,863,73,0,0,53),
    @("Vietnam",964,2,465,475,0,0,24),
    @("Republica del Chad",956,0,865,15,0,0,76),
    @("Lesoto",903,0,271,607,0,0,25),
    @("Santo Tome y Principe",885,0,818,52,0,0,15),
    @("Reunion",855,0,657,193,0,0,5),
    @("Crucero",712,0,651,48,0,0,13),
    @("Guyana",709,0,349,337,0,0,23),
    @("San Marino",699,0,657,0,0,0,42),
    @("Trinidad yTobago",552,0,140,401,0,0,11),
    @("Tanzania",509,0,183,305,0,0,21),
    @("Taiwan",485,1,450,28,0,0,7),
    @("Belice",452,0,35,414,0,0,3),
    @("Guadalupe",446,0,289,143,0,0,14),
    @("Burundi",413,0,315,97,0,0,1),
    @("Comoras",405,0,379,19,0,0,7),
    @("Birmania",375,0,329,40,0,0,6),
    @("Islas Feroe",372,0,225,147,0,0,0),
    @("Mauricio",346,0,334,2,0,0,10),
    @("Martinica",336,0,98,222,0,0,16),
    @("Isla de Man",336,0,312,0,0,0,24),
    @("Papua Nueva Guinea",323,0,110,210,0,0,3),
    @("San Martin (Parte Holandesa)",317,17,107,193,0,0,17),
    @("Islas Turcas y Caicos",298,0,55,241,0,0,2),
    @("Mongolia",298,0,276,22,0,0,0),
    @("Eritrea",285,0,248,37,0,0,0),
    @("Camboya",273,0,238,35,0,0,0),
    @("Gibraltar",215,0,190,25,0,0,0),
    @("Islas Caimanes",203,0,202,0,0,0,1),
    @("Polinesia Francesa",166,0,64,102,0,0,0),
    @("Bermudas",162,0,147,6,0,0,9),
    @("Barbados",151,0,122,22,0,0,7),
    @("Monaco",146,0,114,28,0,0,4),
    @("Brunei",142,0,138,1,0,0,3),
    @("Butan",141,3,103,38,0,0,0),
    @("Seychelles",127,0,126,1,0,0,0),
    @("San Martin (Parte Francesa)",109,0,49,56,0,0,4),
    @("Liechtenstein",94,3,87,6,0,0,1),
    @("Antigua y Barbuda",93,0,88,2,0,0,3),
    @("San Vicente y las Granadinas",57,0,55,2,0,0,0),
    @("Macao",46,0,46,0,0,0,0),
    @("Puerto Rico",39,0,1,36,0,0,2),
    @("Curazao",34,0,31,2,0,0,1),
    @("Guam",32,0,0,31,0,0,1),
    @("Fiyi",28,0,20,7,0,0,1),
    @("Timor Oriental",25,0,24,1,0,0,0),
    @("Santa Lucia",25,0,25,0,0,0,0),
    @("Granada",24,0,23,1,0,0,0),
    @("Nueva Caledonia",23,0,22,1,0,0,0),
    @("Laos",22,0,19,3,0,0,0),
    @("Dominica",18,0,18,0,0,0,0),
    @("Islas Virgenes de los Estados Unidos",17,0,0,17,0,0,0),
    @("San Cristobal y Nieves",17,0,17,0,0,0,0),
    @("Groenlandia",14,0,14,0,0,0,0),
    @("Bonaire, San Eustaquio y Saba",13,0,7,6,0,0,0),
    @("San Bartolome",13,0,9,4,0,0,0),
    @("Montserrat",13,0,12,0,0,0,1),
    @("Islas Malvinas",13,0,13,0,0,0,0),
    @("Santa Sede",12,0,12,0,0,0,0),
    @("Sahara Occidental",10,0,8,1,0,0,1),
    @("Islas Virgenes Britanicas",9,0,7,1,0,0,1),
    @("San Pedro y Miquelon",4,0,1,3,0,0,0),
    @("Anguila",3,0,3,0,0,0,0)
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
}

Write-Output "Updated $($data.Count) rows"
